$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "E3"  = 16.14010000000001
    "E4"  = 16.5466
    "B11" = 5.497300000000001
    "B12" = 5.1904
    "E14" = 16.53680000000001
    "B15" = 5.219999999999996
    "E26" = 16.12219999999999
    "B27" = 5.712699999999999
    "B28" = 5.652799999999997
    "B31" = 4.893799999999999
    "E31" = 16.4069
    "B32" = 6.424199999999997
    "E35" = 16.75369999999999
    "B36" = 8.321500000000007
    "E37" = 16.72130000000001
    "B38" = 4.825899999999994
    "E39" = 16.4244
    "E40" = 17.09760000000001
    "E45" = 16.5082
    "B46" = 6.222500000000004
    "E52" = 16.99240000000001
    "B54" = 4.678999999999998
    "B55" = 5.341699999999996
    "B56" = 5.008299999999998
    "E57" = 16.7337
    "B67" = 5.579099999999996
    "B69" = 5.499499999999995
    "B72" = 5.494700000000002
    "B73" = 8.441699999999994
    "E81" = 16.3159
    "B83" = 5.737799999999996
    "E83" = 16.70979999999999
    "B86" = 5.001200000000001
    "B91" = 5.447700000000001
    "B93" = 5.967200000000002
    "B99" = 4.600799999999999
    "E100" = 16.4442
    "E102" = 16.8349
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
